$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.655, 0.24, -0.008, -105, -22, -100),
    @(0.194, -0.078, -0.541, 97, -36, 20),
    @(-0.544, -0.09, 0.12, 131, 52, 148),
    @(-0.131, -0.006, -0.021, 87, 7, 25),
    @(-0.433, -0.125, 0.757, 16, -29, -67),
    @(-0.554, 0.151, 0.14, 97, 61, -90),
    @(0.193, 0.231, 0.584, -81, -9, 0),
    @(0.339, 0.203, -0.149, -61, 33, 115),
    @(-0.268, -0.156, -0.211, 175, -28, 36),
    @(-0.144, -0.355, -0.435, -126, 15, 106)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
